$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Total" row (row 5) to make room for
# the new team member, pushing Total/disclaimer down by one row.
$ws.Rows.Item(5).Insert()

# Update matriculation/name column and give everyone an equal Code Review 1
# score of 25 (was 33/34/33 before the new member joined).
$ws.Range("A2").Value = "40443267 (Euan Campbell)"
$ws.Range("B2").Value = 25

$ws.Range("A3").Value = "40443517 (Karl Denison)"
$ws.Range("B3").Value = 25

$ws.Range("A4").Value = "40491512 (Joe Black)"
$ws.Range("B4").Value = 25

$ws.Range("A5").Value = "40485296 (Usmaan Chohan)"
$ws.Range("B5").Value = 25

# Move the active selection to match the saved view state.
$ws.Range("E14").Select()
